$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 0.01483866666666667
$ws.Range("N2").Value = 0.044516
$ws.Range("Q2").Value = 0.009467000086222222
$ws.Range("R2").Value = 0.085203000776
